$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "تعارضات الجدول"

# Update header row
$ws.Cells.Item(1, 1).Value = "نوع التعارض"
$ws.Cells.Item(1, 2).Value = "تفاصيل التعارض"

# Widen column B to fit the longer detail text
$ws.Columns.Item(2).ColumnWidth = 239.1666667

# Conflict rows data: column A = conflict type, column B = conflict detail
$data = @(
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ا. خالد الكحسه, اليوم: الخميس, الوقت: 8-10, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الأحد 8-10, الاثنين 10-12, الاثنين 12-2, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: إيهاب الحروي, اليوم: الاثنين, الوقت: 12-2, الأوقات المتاحة: الأربعاء 10-12, الأربعاء 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: هديل دماج, اليوم: الثلاثاء, الوقت: 10-12, الأوقات المتاحة: الخميس 10-12, الخميس 12-2"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: د نشوان المجمر, اليوم: الخميس, الوقت: 12-2, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الأحد 8-10, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: د. فهد الاغبري, اليوم: الأربعاء, الوقت: 12-2, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الأحد 8-10, الاثنين 10-12, الاثنين 12-2, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ا. خالد الكحسه, اليوم: السبت, الوقت: 12-2, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الأحد 8-10, الاثنين 10-12, الاثنين 12-2, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: خولة صادق, اليوم: السبت, الوقت: 10-12, الأوقات المتاحة: الخميس 10-12, الخميس 12-2, السبت 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: عبدالقوي الفقيه, اليوم: الاثنين, الوقت: 10-12, الأوقات المتاحة: الأربعاء 10-12, الأربعاء 12-2, الأربعاء 8-10, الخميس 10-12, الخميس 12-2, الخميس 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ا. معافر منصور, اليوم: الخميس, الوقت: 12-2, الأوقات المتاحة: الأحد 10-12, الأحد 8-10, الاثنين 12-2"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: د. ندى الحميدي, اليوم: الخميس, الوقت: 8-10, الأوقات المتاحة: الثلاثاء 10-12, الثلاثاء 12-2, الثلاثاء 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: إسماعيل الصباحي, اليوم: السبت, الوقت: 12-2, الأوقات المتاحة: الخميس 12-2, الخميس 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: عبده صالح, اليوم: الخميس, الوقت: 8-10, الأوقات المتاحة: الثلاثاء 12-2, الثلاثاء 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: خديجة الشيخ, اليوم: السبت, الوقت: 12-2, الأوقات المتاحة: الاثنين 10-12, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: علي الحمادي, اليوم: الأربعاء, الوقت: 12-2, الأوقات المتاحة: الخميس 10-12, السبت 10-12, السبت 12-2, السبت 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ا. وفاء بركات, اليوم: السبت, الوقت: 12-2, الأوقات المتاحة: الأربعاء 10-12, الأربعاء 8-10, الاثنين 10-12, الاثنين 12-2, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: د. احمد مهيوب, اليوم: الأربعاء, الوقت: 8-10, الأوقات المتاحة: الاثنين 10-12, الاثنين 12-2, الاثنين 8-10, الثلاثاء 8-10, الخميس 10-12, الخميس 12-2"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ليلى الشعيبي, اليوم: الخميس, الوقت: 10-12, الأوقات المتاحة: الأربعاء 10-12, الأربعاء 12-2, الأربعاء 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: د. محمد خشافه, اليوم: الأربعاء, الوقت: 10-12, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الاثنين 10-12, الاثنين 12-2, الاثنين 8-10, الثلاثاء 10-12, الثلاثاء 12-2, الثلاثاء 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: د. ناجي غلاب, اليوم: الأربعاء, الوقت: 12-2, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الأحد 8-10, الاثنين 10-12, الاثنين 12-2, الاثنين 8-10, الثلاثاء 10-12, الثلاثاء 12-2, الثلاثاء 8-10, الخميس 10-12, الخميس 12-2, الخميس 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: محمد الطيب, اليوم: الخميس, الوقت: 10-12, الأوقات المتاحة: الأربعاء 8-10, الثلاثاء 12-2"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ا. وفاء بركات, اليوم: الأحد, الوقت: 10-12, الأوقات المتاحة: الأربعاء 10-12, الأربعاء 8-10, الاثنين 10-12, الاثنين 12-2, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: لينا الفقية, اليوم: الاثنين, الوقت: 12-2, الأوقات المتاحة: الاثنين 10-12, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: د. عبدالباسط الادريسي, اليوم: الثلاثاء, الوقت: 10-12, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الاثنين 10-12"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ا. فؤاد الهمداني, اليوم: الخميس, الوقت: 12-2, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الأحد 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ا.هناء العبيدي, اليوم: الخميس, الوقت: 10-12, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الأحد 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: د. عادل العفيري, اليوم: الثلاثاء, الوقت: 12-2, الأوقات المتاحة: الاثنين 10-12, الاثنين 12-2, الاثنين 8-10"),
    @("تم تعيين دكتور في غير موعدة", "اسم الدكتور: ا. صلاح امين, اليوم: الخميس, الوقت: 10-12, الأوقات المتاحة: الأحد 10-12, الأحد 12-2, الأحد 8-10, الاثنين 10-12, الاثنين 8-10")

)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
